$d = $word.ActiveDocument

# Helper: locate the 1-based index of the first paragraph (at/after $startIdx)
# whose style is $styleName and whose text contains $textSnippet.
function Find-ParaIndex($startIdx, $styleName, $textSnippet) {
    for ($i = $startIdx; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Style.NameLocal -eq $styleName -and $p.Range.Text.Contains($textSnippet)) {
            return $i
        }
    }
    return -1
}

# --- Change 1: heading "07/04/22 - 11/04/22" -> "07/04/22 - 12/04/22" ---
$idxHeading1 = Find-ParaIndex 1 "Heading 1" "11/04/22"
$pHeading1 = $d.Paragraphs.Item($idxHeading1)
$pHeading1.Range.Find.Execute("11/04/22", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "12/04/22", 2)

# --- Change 2: the blank "- " note right after it becomes the "Slow reading..." entry ---
$pNote1 = $d.Paragraphs.Item($idxHeading1 + 1)
$rNote1 = $d.Range($pNote1.Range.Start, $pNote1.Range.End - 1)
$rNote1.Text = "Slow reading of the S2S paper ; global understanding of the proof of completeness, still need to go in detail in the construction of partial tableaux and why they fit."

# --- Change 3: the next "06/04/22" heading becomes "13/04/22" ---
$idxHeading2 = Find-ParaIndex ($idxHeading1 + 1) "Heading 1" "06/04/22"
$pHeading2 = $d.Paragraphs.Item($idxHeading2)
$rHeading2 = $d.Range($pHeading2.Range.Start, $pHeading2.Range.End - 1)
$rHeading2.Text = "13/04/22"

# --- Change 4: the blank "- " note right after it gets real content ---
$pNote2 = $d.Paragraphs.Item($idxHeading2 + 1)
$rNote2 = $d.Range($pNote2.Range.Start, $pNote2.Range.End - 1)
$rNote2.Text = "- Accessed the full version of the S2S paper (with the proofs)."

# --- Change 5: insert a brand-new paragraph with a second note right after it ---
$pNote2.Range.InsertParagraphAfter()
$pNote3 = $d.Paragraphs.Item($idxHeading2 + 2)
$rNote3 = $d.Range($pNote3.Range.Start, $pNote3.Range.End - 1)
$rNote3.Text = "- Back to some level of productivity. On Sam’s advice, rewriting the paper in my own way, to explicit the proofs, to get used with LaTeX, and most importantly to have a deep understanding of the content."
